# Saldo.xlsx update — "Add files via upload"
# Replaces the two smallest-id "top" accounts with a fresh trio of larger
# balances, inserts a new GUSTAVO row up near the top of the list, and
# removes the two stale rows (the old DILSON and old GUSTAVO entries)
# whose balances migrated to the new rows above.
#
# Account numbers in column A are zero-padded digit strings, so they are
# written with a leading apostrophe to force Excel to keep them as text
# (otherwise the leading zeros would be silently stripped by numeric
# auto-conversion), matching the workbook's existing inlineStr cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 3 (was 008007764 / LUIS / 40999.9) -> 004472404 / DILSON / 56773.27
$ws.Range("A3").Value = "'004472404"
$ws.Range("B3").Value = "DILSON"
$ws.Range("C3").Value = 56773.27

# 2) Row 4 (was 008004799 / OLGA / 40000) -> 003301389 / EDMUNDO / 54727.43
$ws.Range("A4").Value = "'003301389"
$ws.Range("B4").Value = "EDMUNDO"
$ws.Range("C4").Value = 54727.43

# 3) Insert a brand-new row 5: 004363260 / LARISSA / 23234.07
$ws.Rows(5).Insert()
$ws.Range("A5").Value = "'004363260"
$ws.Range("B5").Value = "LARISSA"
$ws.Range("C5").Value = 23234.07

# 4) Insert a brand-new row 9 (just above THAYSA): 004565108 / GUSTAVO / 1976.91
$ws.Rows(9).Insert()
$ws.Range("A9").Value = "'004565108"
$ws.Range("B9").Value = "GUSTAVO"
$ws.Range("C9").Value = 1976.91

# 5) Remove the old DILSON row (004472404 / DILSON / 805.17), now at row 21
$ws.Rows(21).Delete()

# 6) Remove the old GUSTAVO row (004565108 / GUSTAVO / 476.91), now at row 38
$ws.Rows(38).Delete()
